$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new value "x" in cell C6 (next to B6's "fix finish message..." row)
$ws.Range("C6").Value = "x"

# Move the active selection to F26 to match the diff's sheetView selection
$ws.Range("F26").Select()
